# Computer Management | 添加 "端口" (Port) worksheet, tidy up "依赖性" sheet selection.
$wb = $excel.ActiveWorkbook

$wsExit = $wb.Worksheets.Item(1)     # "Exit Code" - style donors live here
$wsDep  = $wb.Worksheets.Item(2)     # "依赖性"

# ------------------------------------------------------------------
# Add the new "端口" sheet as the last tab and make it the active one
# (mirrors the workbook.xml <sheets>/activeTab change in the diff).
# ------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsPort = $wb.Worksheets.Add($null, $lastSheet)
$wsPort.Name = "端口"

$wsPort.Columns.Item(1).ColumnWidth = 17

# -- A1 : section header, reuse the shaded "fillId2" look used
#    throughout "Exit Code" (e.g. A3 there) so the style is deduped
#    against the existing xf instead of minting a new fill.
$wsExit.Range("A3").Copy() | Out-Null
$wsPort.Range("A1").PasteSpecial(-4122) | Out-Null
$wsPort.Range("A1").Value = "Computer"

# -- A2 / B2 : ComputerUI port row (string must land in the shared
#    string table before "Control" to match authoring order)
$wsPort.Range("A2").Value = "ComputerUI"
$wsPort.Range("B2").Value = 23332

# -- A4 : second section header, same shaded style as A1
$wsExit.Range("A3").Copy() | Out-Null
$wsPort.Range("A4").PasteSpecial(-4122) | Out-Null
$wsPort.Range("A4").Value = "Control"

# -- B2 / A3 : centered cells, reuse the existing plain "center" xf
#    already present in the workbook (依赖性!B3 carries it).
$wsDep.Range("B3").Copy() | Out-Null
$wsPort.Range("B2").PasteSpecial(-4122) | Out-Null
$wsPort.Range("B2").Value = 23332

$wsDep.Range("B3").Copy() | Out-Null
$wsPort.Range("A3").PasteSpecial(-4122) | Out-Null

# -- A5 : "用于进程间通讯" label, left + vertically centered (new style)
$wsPort.Range("A5").Value = "用于进程间通讯"
$wsPort.Range("A5").HorizontalAlignment = -4131
$wsPort.Range("A5").VerticalAlignment = -4108

# -- B5 : port number, centered + vertically centered (new style)
$wsPort.Range("B5").Value = 23333
$wsPort.Range("B5").HorizontalAlignment = -4108
$wsPort.Range("B5").VerticalAlignment = -4108

$wsPort.Range("B5").Select() | Out-Null

# ------------------------------------------------------------------
# "依赖性" keeps its layout; it just stops being the selected tab
# once "端口" becomes the active sheet (tabSelected moves off it).
# ------------------------------------------------------------------
$wsDep.Range("F9").Select() | Out-Null

$wsPort.Activate() | Out-Null
